$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'" + '34.506.95'
$ws.Range("E2").Value = '  -0.19%  '
$ws.Range("D3").Value = "'" + '1.807.20'
$ws.Range("E3").Value = '  -0.40%  '
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").Value = "'" + '225.61'
$ws.Range("E5").Value = '  -1.53%  '
$ws.Range("E6").Value = '  +2.54%  '
$ws.Range("E7").Value = '  +0.06%  '
$ws.Range("D8").Value = "'" + '37.37'
$ws.Range("E8").Value = '  +6.85%  '
$ws.Range("D9").Value = "'" + '0.291'
$ws.Range("E9").Value = '  -3.57%  '
$ws.Range("D10").Value = "'" + '0.0680'
$ws.Range("E10").Value = '  -2.54%  '
$ws.Range("E11").Value = '  +1.41%  '
$ws.Range("D12").Value = "'" + '2.067.84'
$ws.Range("E12").Value = '  -0.42%  '
$ws.Range("D13").Value = "'" + '11.29'
$ws.Range("E13").Value = '  +0.42%  '
$ws.Range("D14").Value = "'" + '1.813.33'
$ws.Range("E14").Value = '  -0.04%  '
$ws.Range("D15").Value = "'" + '0.632'
$ws.Range("D16").Value = "'" + '34.457.69'
$ws.Range("E16").Value = '  -0.29%  '
$ws.Range("D17").Value = "'" + '4.42'
$ws.Range("E17").Value = '  -0.62%  '
$ws.Range("D18").Value = "'" + '68.56'
$ws.Range("E18").Value = '  -1.03%  '
$ws.Range("D19").Value = "'" + '243.51'
$ws.Range("E19").Value = '  -1.08%  '
$ws.Range("E20").Value = '  -3.09%  '
$ws.Range("D21").Value = "'" + '11.20'
$ws.Range("E21").Value = '  -2.56%  '
$ws.Range("E22").Value = '  +0.01%  '
$ws.Range("D23").Value = "'" + '4.13'
$ws.Range("E23").Value = '  -1.04%  '
$ws.Range("E24").Value = '  +3.88%  '
$ws.Range("E25").Value = '  -1.13%  '
$ws.Range("D26").Value = "'" + '7.84'
$ws.Range("E26").Value = '  -1.22%  '
$ws.Range("D27").Value = "'" + '17.31'
$ws.Range("E27").Value = '  +2.67%  '
$ws.Range("E28").Value = '  +0.46%  '
$ws.Range("E29").Value = '  +0.08%  '
$ws.Range("D30").Value = "'" + '3.81'
$ws.Range("E30").Value = '  -1.16%  '
$ws.Range("E33").Value = '  -2.83%  '
$ws.Range("E34").Value = '  -2.07%  '
$ws.Range("D35").Value = "'" + '1.362.55'
$ws.Range("E35").Value = '  -2.54%  '
$ws.Range("D36").Value = "'" + '0.654'
$ws.Range("E36").Value = '  -4.76%  '
$ws.Range("E37").Value = '  +0.72%  '
$ws.Range("D38").Value = "'" + '2.36'
$ws.Range("E38").Value = '  -5.86%  '
$ws.Range("E39").Value = '  -2.21%  '
$ws.Range("D40").Value = "'" + '2.43'
$ws.Range("E40").Value = '  +1.77%  '
$ws.Range("E41").Value = '  -2.27%  '
$ws.Range("D42").Value = "'" + '80.61'
$ws.Range("E42").Value = '  -3.86%  '
$ws.Range("D43").Value = "'" + '0.937'
$ws.Range("E43").Value = '  -2.70%  '
$ws.Range("D44").Value = "'" + '1.19'
$ws.Range("E44").Value = '  +6.17%  '
$ws.Range("D45").Value = "'" + '13.72'
$ws.Range("E45").Value = '  +0.99%  '
$ws.Range("D46").Value = "'" + '0.0502'
$ws.Range("E46").Value = '  -2.24%  '
$ws.Range("D47").Value = "'" + '1.968.09'
$ws.Range("E47").Value = '  -0.38%  '
$ws.Range("D48").Value = "'" + '5.80'
$ws.Range("E48").Value = '  -3.27%  '
$ws.Range("E49").Value = '  +0.05%  '
$ws.Range("D50").Value = "'" + '102.56'
$ws.Range("E50").Value = '  -2.39%  '
$ws.Range("E51").Value = '  -6.30%  '

# Row 31/32 swap (InternetComputer(DFINITY) <-> PancakeSwap)
$ws.Range("B31").Value = 'PancakeSwap'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D31").Value = "'" + '1.23'
$ws.Range("E31").Value = '  -1.21%  '

$ws.Range("B32").Value = 'InternetComputer(DFINITY)'
$ws.Range("C32").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D32").Value = "'" + '3.91'
$ws.Range("E32").Value = '  -3.32%  '
